# Add three new Star Wars characters (R2-D2, Grogu, Darth Maul) to the
# character table as rows 5, 6 and 7, matching the formatting already
# used for the existing rows (2-4) and wiring up the Wiki Link / Image
# Link hyperlinks for each new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting (styles) from row 4 down onto the three new
#     rows first, so the new cells pick up the same "body" / "hyperlink"
#     cell styles used elsewhere in the sheet before any values or
#     hyperlinks are added.
$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A5:H7").PasteSpecial(-4122) | Out-Null

# --- Row 5: R2-D2 ---------------------------------------------------
$ws.Range("A5").Value = "R2-D2"
$ws.Range("B5").Value = "Droid"
$ws.Range("C5").Value = "Astromech"
$ws.Range("D5").Value = "~32 BBY"
$ws.Range("E5").Value = "Naboo"
$ws.Range("F5").Value = "Episode IV: A New Hope (1977)"
$ws.Range("G5").Value = "https://starwars.fandom.com/wiki/R2-D2"
$ws.Range("H5").Value = "https://pyxis.nymag.com/v1/imgs/7ef/846/3245bc42a87b290d806985f75dc6e7bd4a-28-r2-d2.rsquare.w330.jpg"

# --- Row 6: Grogu (The Child) ---------------------------------------
$ws.Range("A6").Value = "Grogu (The Child)"
$ws.Range("B6").Value = "Unknown"
$ws.Range("C6").Value = "Male"
$ws.Range("D6").Value = "41 BBY"
$ws.Range("E6").Value = "Unknown"
$ws.Range("F6").Value = "The Mandalorian (Season 1, 2019)"
$ws.Range("G6").Value = "https://starwars.fandom.com/wiki/Grogu"
$ws.Range("H6").Value = "https://pyxis.nymag.com/v1/imgs/99e/6f5/6eed622d1b1b0a77caad3e658d61630b76-baby-yoda.rsquare.w700.jpg"

# --- Row 7: Darth Maul -----------------------------------------------
$ws.Range("A7").Value = "Darth Maul"
$ws.Range("B7").Value = "Zabrak"
$ws.Range("C7").Value = "Male"
$ws.Range("D7").Value = "54 BBY"
$ws.Range("E7").Value = "Dathomir"
$ws.Range("F7").Value = "Epsiode I: The Phantom Menace (1999)"
$ws.Range("G7").Value = "https://starwars.fandom.com/wiki/Maul"
$ws.Range("H7").Value = "https://i.pinimg.com/736x/cc/bb/02/ccbb029fca3b6bf6256b49acc110a200.jpg"

# --- Hyperlinks for the Wiki Link (G) / Image Link (H) columns -------
$ws.Hyperlinks.Add($ws.Range("G5"), "https://starwars.fandom.com/wiki/R2-D2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H5"), "https://pyxis.nymag.com/v1/imgs/7ef/846/3245bc42a87b290d806985f75dc6e7bd4a-28-r2-d2.rsquare.w330.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G6"), "https://starwars.fandom.com/wiki/Grogu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H6"), "https://pyxis.nymag.com/v1/imgs/99e/6f5/6eed622d1b1b0a77caad3e658d61630b76-baby-yoda.rsquare.w700.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G7"), "https://starwars.fandom.com/wiki/Maul") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H7"), "https://i.pinimg.com/736x/cc/bb/02/ccbb029fca3b6bf6256b49acc110a200.jpg") | Out-Null

# Adding a hyperlink re-applies the builtin "Hyperlink" style (blue,
# underlined) to the cell, which would overwrite the plain centered
# style used by the rest of the table - paste the correct formatting
# back on top of the Wiki Link / Image Link cells to restore it. Each
# destination is pasted from a single source cell so the paste can't
# spill into the next column.
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null

# --- Widen the "First Screen Appearance" / "Wiki Link" columns a touch
$ws.Columns("F").ColumnWidth = 36.3
$ws.Columns("G").ColumnWidth = 50.3

# --- Move the active selection down below the new data ---------------
$ws.Range("A8").Select() | Out-Null
